# Update LR-pairs worksheet (Tgfb2-Eng) with refreshed TPM-based NATMI output.
# A new "Resolving-Mac" sending/target cluster is introduced, extending the
# data table from 12 to 16 data rows (plus header), and several numeric
# metrics are recalculated with the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate,
# G=Ligand average expression value, H=Ligand total expression value,
# I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate,
# M=Receptor average expression value, N=Receptor total expression value,
# O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity
$rows = @(
    @("ECs", "Tgfb2", "Eng", "ECs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 218.721583, 656.164749, 0.7793342808141792, 0.7793342808141792, 560.8356761213934, 5047.52108509254, 0.04090842183931226, 0.04090842183931226),
    @("ECs", "Tgfb2", "Eng", "FAPs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 46.29469433333333, 138.884083, 0.164954193449581, 0.164954193449581, 118.7066947904645, 1068.36025311418, 0.00865869229151482, 0.00865869229151482),
    @("ECs", "Tgfb2", "Eng", "MuSCs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 8.515309999999999, 25.54593, 0.03034118948727519, 0.03034118948727519, 21.83456052086667, 196.5110446878, 0.001592654409293088, 0.001592654409293087),
    @("ECs", "Tgfb2", "Eng", "Resolving-Mac", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 7.120231, 21.360693, 0.02537033624896462, 0.02537033624896462, 18.25736405275333, 164.31627647478, 0.001331726889254218, 0.001331726889254218),
    @("FAPs", "Tgfb2", "Eng", "ECs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 218.721583, 656.164749, 0.7793342808141792, 0.7793342808141792, 3629.643696041785, 32666.79326437606, 0.2647531206840318, 0.2647531206840317),
    @("FAPs", "Tgfb2", "Eng", "FAPs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 46.29469433333333, 138.884083, 0.164954193449581, 0.164954193449581, 768.2517799222616, 6914.266019300353, 0.05603774729384325, 0.05603774729384323),
    @("FAPs", "Tgfb2", "Eng", "MuSCs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 8.515309999999999, 25.54593, 0.03034118948727519, 0.03034118948727519, 141.3099742485933, 1271.78976823734, 0.01030741852344742, 0.01030741852344742),
    @("FAPs", "Tgfb2", "Eng", "Resolving-Mac", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 7.120231, 21.360693, 0.02537033624896462, 0.02537033624896462, 118.1588995883927, 1063.430096295534, 0.008618735066676909, 0.008618735066676908),
    @("MuSCs", "Tgfb2", "Eng", "ECs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 218.721583, 656.164749, 0.7793342808141792, 0.7793342808141792, 6488.258962369679, 58394.33066132711, 0.4732659599527079, 0.4732659599527079),
    @("MuSCs", "Tgfb2", "Eng", "FAPs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 46.29469433333333, 138.884083, 0.164954193449581, 0.164954193449581, 1373.307386031559, 12359.76647428403, 0.1001716550048112, 0.1001716550048112),
    @("MuSCs", "Tgfb2", "Eng", "MuSCs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 8.515309999999999, 25.54593, 0.03034118948727519, 0.03034118948727519, 252.60212397446, 2273.41911577014, 0.01842527978340799, 0.01842527978340799),
    @("MuSCs", "Tgfb2", "Eng", "Resolving-Mac", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 7.120231, 21.360693, 0.02537033624896462, 0.02537033624896462, 211.217850411646, 1900.960653704814, 0.01540663208943596, 0.01540663208943595),
    @("Resolving-Mac", "Tgfb2", "Eng", "ECs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 218.721583, 656.164749, 0.7793342808141792, 0.7793342808141792, 5.576744201751, 50.190697815759, 0.0004067783381273135, 0.0004067783381273135),
    @("Resolving-Mac", "Tgfb2", "Eng", "FAPs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 46.29469433333333, 138.884083, 0.164954193449581, 0.164954193449581, 1.180375821417, 10.623382392753, 0.00008609885941171745, 0.00008609885941171745),
    @("Resolving-Mac", "Tgfb2", "Eng", "MuSCs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 8.515309999999999, 25.54593, 0.03034118948727519, 0.03034118948727519, 0.21711485907, 1.95403373163, 0.00001583677112669258, 0.00001583677112669257),
    @("Resolving-Mac", "Tgfb2", "Eng", "Resolving-Mac", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 7.120231, 21.360693, 0.02537033624896462, 0.02537033624896462, 0.181544529807, 1.633900768263, 0.00001324220359754154, 0.00001324220359754153)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $rows[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
